$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (Ахматова Светлана) - ДЗ_1..ДЗ_3 graded
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 5

# Row 15 (Петров Артём) - ДЗ_1..ДЗ_3 graded
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = 5

# Row 22 (Соди Гитанджелина) - ДЗ_1..ДЗ_2 graded
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 5

# Move active selection to F15 to match the saved view state
$ws.Range("F15").Select()
